$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I18").Value = -0.6067089621107076
$ws.Range("J18").Value = 0.1224968750237004
$ws.Range("K18").Value = 0.006238849779163275
$ws.Range("L18").Value = 2.054960549259019
